# Update cryptocurrency price/volume data per the Oct 30 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "226.47") must be forced to stay
# text (the source sheet stores every Price/Volume cell as a string), so we temporarily
# mark them as Text format, assign the literal string, then clear the temporary format so
# the cell keeps its original (default) style once the text value has "stuck".
$textForced = @(
    "D5",
    "D8",
    "D10",
    "D13",
    "D15",
    "D18",
    "D19",
    "D21",
    "D24",
    "D25",
    "D26",
    "D27",
    "D30",
    "D31",
    "D32",
    "D33",
    "D37",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D49",
)
foreach ($addr in $textForced) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row-by-row value updates (Coin / Link / Price / Volume(1h))
$ws.Range("D2").Value = '34.354.62'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '1.786.84'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '226.47'
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("E6").Value = '  +2.12%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '32.94'
$ws.Range("E8").Value = '  +3.64%  '

$ws.Range("E9").Value = '  +1.15%  '

$ws.Range("D10").Value = '0.0688'
$ws.Range("E10").Value = '  +0.35%  '

$ws.Range("D12").Value = '2.044.82'
$ws.Range("E12").Value = '  +0.40%  '

$ws.Range("D13").Value = '11.21'
$ws.Range("E13").Value = '  +2.71%  '

$ws.Range("D14").Value = '1.794.05'
$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("D15").Value = '0.635'
$ws.Range("E15").Value = '  +2.11%  '

$ws.Range("D16").Value = '34.322.35'
$ws.Range("E16").Value = '  +0.69%  '

$ws.Range("E17").Value = '  +2.62%  '

$ws.Range("D18").Value = '68.38'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").Value = '245.01'
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").Value = '0.0₃0794'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("D21").Value = '11.27'
$ws.Range("E21").Value = '  +3.45%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D24").Value = '169.07'
$ws.Range("E24").Value = '  +4.65%  '

$ws.Range("D25").Value = '2.05'
$ws.Range("E25").Value = '  +0.83%  '

$ws.Range("D26").Value = '7.32'
$ws.Range("E26").Value = '  +3.15%  '

$ws.Range("D27").Value = '16.53'
$ws.Range("E27").Value = '  +1.88%  '

$ws.Range("E28").Value = '  +1.75%  '

$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").Value = '4.04'
$ws.Range("E30").Value = '  +9.31%  '

$ws.Range("D31").Value = '0.0526'
$ws.Range("E31").Value = '  +1.98%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '3.79'
$ws.Range("E33").Value = '  +2.68%  '

$ws.Range("E34").Value = '  +1.33%  '

$ws.Range("D35").Value = '1.410.28'
$ws.Range("E35").Value = '  -2.52%  '

$ws.Range("E36").Value = '  +5.12%  '

$ws.Range("D37").Value = '0.683'
$ws.Range("E37").Value = '  +5.04%  '

$ws.Range("E38").Value = '  +3.15%  '

$ws.Range("D39").Value = '0.0191'
$ws.Range("E39").Value = '  +0.39%  '

$ws.Range("D40").Value = '84.42'
$ws.Range("E40").Value = '  +5.19%  '

$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("D42").Value = '2.78'
$ws.Range("E42").Value = '  +2.49%  '

$ws.Range("D43").Value = '0.939'
$ws.Range("E43").Value = '  +2.75%  '

$ws.Range("D44").Value = '13.99'

$ws.Range("D45").Value = '0.0528'
$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("E46").Value = '  +2.63%  '

$ws.Range("D47").Value = '6.08'
$ws.Range("E47").Value = '  +0.64%  '

$ws.Range("D48").Value = '1.945.32'
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("D49").Value = '105.36'
$ws.Range("E49").Value = '  +1.11%  '

$ws.Range("E50").Value = '  -0.09%  '

$ws.Range("E51").Value = '  -1.61%  '

# Drop the temporary Text format now that the literal strings are committed, restoring
# each cell to its original (default/general) style.
foreach ($addr in $textForced) {
    $ws.Range($addr).ClearFormats()
}
